# Updated symbol list on Tue Dec 27 07:31:23 UTC 2022 with GitHub Actions
# Applies updated prices / labels / swapped rows for KickToken <-> BKEXToken
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve each cell's original Text storage type (these columns are inline
# strings in the workbook, not numbers) by forcing Text format before the
# write, then stripping the format again so no stray style id is left behind.
function Set-TextValue($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = '@'
    $rng.Value = $val
    $rng.ClearFormats()
}

Set-TextValue 'D2' '243.49'
Set-TextValue 'D4' '5.403'
Set-TextValue 'D7' '6.498'
Set-TextValue 'D8' '0.8123'
Set-TextValue 'D9' '0.9261'
Set-TextValue 'D11' '0.07446'
Set-TextValue 'D12' '0.03321'
Set-TextValue 'D13' '0.03063'
Set-TextValue 'D15' '3.860'
Set-TextValue 'D16' '0.001572'
Set-TextValue 'D18' '0.0005888'
Set-TextValue 'E18' '17OneONE'
Set-TextValue 'D19' '0.005897'
Set-TextValue 'E20' '19BitKanKANBestin24h'
Set-TextValue 'D21' '0.004880'
Set-TextValue 'D22' '0.00007999'
Set-TextValue 'E22' '21NitroExNTX'
Set-TextValue 'D23' '3.572'
Set-TextValue 'D26' '0.1308'
Set-TextValue 'D27' '0.0002338'
Set-TextValue 'D40' '0.03947'
Set-TextValue 'B41' 'BKEXToken'
Set-TextValue 'C41' 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
Set-TextValue 'D41' '0.1078'
Set-TextValue 'E41' '40BKEXTokenBKK'
Set-TextValue 'D42' '0.002660'
Set-TextValue 'B43' 'KickToken'
Set-TextValue 'C43' 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
Set-TextValue 'D43' '0.003066'
Set-TextValue 'E43' '42KickTokenKICKWorstin24h'
Set-TextValue 'D44' '0.008566'
Set-TextValue 'D45' '0.00005198'
Set-TextValue 'D47' '0.6697'
Set-TextValue 'D48' '0.002272'
Set-TextValue 'D49' '0.00002099'
Set-TextValue 'D50' '0.0001999'
